# Update "想去人数" (interest counts) figures re-scraped at a later time.
# Sheet 1 = 展览, Sheet 2 = 演出, Sheet 3 = 本地生活, Sheet 4 = 全部类型

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F8").Value = 196
$ws1.Range("F9").Value = 374
$ws1.Range("F11").Value = 491
$ws1.Range("F12").Value = 525
$ws1.Range("F14").Value = 12319
$ws1.Range("F15").Value = 81
$ws1.Range("F16").Value = 5468

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 115
$ws2.Range("F4").Value = 9

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 115
$ws4.Range("F10").Value = 196
$ws4.Range("F11").Value = 375
$ws4.Range("F13").Value = 491
$ws4.Range("F14").Value = 525
$ws4.Range("F16").Value = 12319
$ws4.Range("F17").Value = 9
$ws4.Range("F18").Value = 81
$ws4.Range("F19").Value = 5468
